$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (row, new text), 1-based row index
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "10845"
$t.Cell(6, 1).Range.Text  = "0.01927"
$t.Cell(7, 1).Range.Text  = "0.00590"
$t.Cell(9, 1).Range.Text  = "0.01927"
$t.Cell(10, 1).Range.Text = "0.01927"
$t.Cell(11, 1).Range.Text = "0.01927"
$t.Cell(12, 1).Range.Text = "1.95856"

# Rows 44-46 previously held a tab-separated run of raw perf-log values;
# collapse each back down to its single summarized figure.
$t.Cell(44, 1).Range.Text = "99.95"
$t.Cell(45, 1).Range.Text = "1.96"
$t.Cell(46, 1).Range.Text = "3891"
